$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $ws.Range("ZZ1").NumberFormat = "@"
    $ws.Range("ZZ1").Value = $text
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("ZZ1").Clear()
}

$ws.Range("D2").Value = "51.597.33"
$ws.Range("D3").Value = "2.988.52"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "381.58"
$ws.Range("E5").Value = "  +2.15%  "
Set-TextValue "D6" "104.29"
$ws.Range("E6").Value = "  +4.15%  "
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue "D9" "0.596"
$ws.Range("E9").Value = "  +2.01%  "
Set-TextValue "D10" "36.76"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").Value = "3.465.20"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D14" "18.52"
$ws.Range("E14").Value = "  +3.17%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D15" "7.84"
$ws.Range("E15").Value = "  +4.16%  "
$ws.Range("D16").Value = "2.985.74"
$ws.Range("E16").Value = "  +1.77%  "
Set-TextValue "D17" "11.22"
$ws.Range("E17").Value = "  +2.60%  "
Set-TextValue "D18" "0.996"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "51.605.33"
$ws.Range("E19").Value = "  +1.53%  "
Set-TextValue "D20" "3.08"
$ws.Range("E20").Value = "  +1.28%  "
Set-TextValue "D21" "12.61"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("E22").Value = "  +1.44%  "
Set-TextValue "D23" "70.51"
$ws.Range("E23").Value = "  +2.32%  "
Set-TextValue "D24" "267.56"
$ws.Range("E24").Value = "  +1.03%  "
Set-TextValue "D25" "3.22"
$ws.Range("E25").Value = "  +2.67%  "
Set-TextValue "D26" "7.90"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("E27").Value = "  +4.88%  "
Set-TextValue "D28" "7.21"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("E29").Value = "  -0.09%  "
Set-TextValue "D30" "26.15"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("E31").Value = "  +1.39%  "
Set-TextValue "D32" "10.39"
$ws.Range("E32").Value = "  +4.89%  "
Set-TextValue "D33" "34.57"
$ws.Range("E33").Value = "  +4.70%  "
Set-TextValue "D34" "51.36"
$ws.Range("E34").Value = "  +1.29%  "
Set-TextValue "D35" "2.06"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("E37").Value = "  +0.03%  "
Set-TextValue "D38" "3.30"
$ws.Range("E38").Value = "  +6.61%  "
Set-TextValue "D39" "16.96"
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("E40").Value = "  +6.44%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D41" "0.117"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "1.84"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D43" "3.86"
$ws.Range("E43").Value = "  +14.95%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D44" "127.48"
$ws.Range("E44").Value = "  +6.68%  "
Set-TextValue "D45" "21.36"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D47" "2.35"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D48" "0.269"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "2.036.36"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "3.287.73"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("E51").Value = "  +2.18%  "
